$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 corresponds to the
# 5d9fc714-b905-471c-9b91-4bbcdf2912ca...zh-cn.xlf entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-21 18:44:01"
$wsZhCn.Range("H4").Value = "2016-03-21 18:44:24"

# de-de sheet: row 4 corresponds to the
# 5d9fc714-b905-471c-9b91-4bbcdf2912ca...de-de.xlf entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-21 18:44:06"
$wsDeDe.Range("H4").Value = "2016-03-21 18:44:30"
